$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Medieval Blacksmith" row (row 15) following the existing
# Name / Image / Link / Price column layout used by the rest of the sheet.
$ws.Range("A15").Value = "Medieval Blacksmith"
$ws.Range("B15").Value = "https://www.lego.com/cdn/cs/set/assets/blt64027ee5a1724b6d/21325_alt1.jpg?format=webply&fit=bounds&quality=75&width=1200&height=1200&dpr=1"
$ws.Range("C15").Value = "https://www.lego.com/en-ch/product/medieval-blacksmith-21325"
$ws.Range("D15").Value = "189 CHF"

# Update the active selection to match the saved view state.
[void]$ws.Range("A17").Select()
